$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 2 (the data row) with plain text values.
$ws.Range("A2").Value = "Yo"
$ws.Range("B2").Value = "is"
$ws.Range("C2").Value = "is"
$ws.Range("D2").Value = "is"
$ws.Range("E2").Value = "is"
$ws.Range("F2").Value = "is"
$ws.Range("G2").Value = "is"
$ws.Range("H2").Value = "isi"
$ws.Range("I2").Value = "is"
$ws.Range("J2").Value = "is"

# Move the active selection to I6 (was C9).
$ws.Range("I6").Select()
